# AlphaFiberF-HW20: rerun simulation -- two new data rows ("Holden" and
# "Rizzie Spiral") inserted right after the "Spiral5" row, all subsequent
# data rows pushed down by two, and one label corrected ("Thomas Hex" ->
# "Matthies Hex").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push existing data rows 4..29 down to 6..31 (columns A:T), working
#    from the bottom up so we never clobber a row before it's been read.
for ($r = 29; $r -ge 4; $r--) {
    $dest = $r + 2
    for ($c = 1; $c -le 20; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dest, $c)
        $dstCell.Value2 = $srcCell.Value2
    }
}

# 2) New row 4 ("Holden") -- A4/B4 already correct (2 / index "2") from the
#    shift above since old row4 also had A=2; only the label text and the
#    measured values need to change.
$ws.Range("B4").Value2 = "Holden"
$ws.Range("C4").Value2 = 0.9910260086134355
$ws.Range("D4").Value2 = 1.025959026844193
$ws.Range("E4").Value2 = 1.014824159203408
$ws.Range("F4").Value2 = 0.9556510368664883
$ws.Range("G4").Value2 = 0.9556510368664883
$ws.Range("H4").Value2 = 1.050326151450064
$ws.Range("I4").Value2 = 1.050326151450064
$ws.Range("J4").Value2 = 1.005749914674527
$ws.Range("K4").Value2 = 0.9556510368664883
$ws.Range("L4").Value2 = 1.005749914674527
$ws.Range("M4").Value2 = 1.028038033062296
$ws.Range("N4").Value2 = 1.028038033062296
$ws.Range("O4").Value2 = 1.023633408442666
$ws.Range("P4").Value2 = 1.00390903433036
$ws.Range("Q4").Value2 = 1.00390903433036
$ws.Range("R4").Value2 = 0.9918445349643921
$ws.Range("S4").Value2 = 0.9918445349643921
$ws.Range("T4").Value2 = 1.007256049608686

# 3) New row 5 ("Rizzie Spiral")
$ws.Range("B5").Value2 = "Rizzie Spiral"
$ws.Range("C5").Value2 = 1.108665750252519
$ws.Range("D5").Value2 = 1.175003648505654
$ws.Range("E5").Value2 = 0.7650352361579889
$ws.Range("F5").Value2 = 1.211960774177592
$ws.Range("G5").Value2 = 1.211960774177592
$ws.Range("H5").Value2 = 0.5486163265634254
$ws.Range("I5").Value2 = 0.5486163265634254
$ws.Range("J5").Value2 = 1.066641852591811
$ws.Range("K5").Value2 = 1.211960774177592
$ws.Range("L5").Value2 = 1.066641852591811
$ws.Range("M5").Value2 = 0.807629089577618
$ws.Range("N5").Value2 = 0.807629089577618
$ws.Range("O5").Value2 = 0.7934311384377416
$ws.Range("P5").Value2 = 0.9424063177776092
$ws.Range("Q5").Value2 = 0.9424063177776092
$ws.Range("R5").Value2 = 1.009794931877605
$ws.Range("S5").Value2 = 1.009794931877605
$ws.Range("T5").Value2 = 0.9793205980414982

# 4) Rename the "Thomas Hex" label (now on row 11, after the two-row
#    shift) to "Matthies Hex".
for ($r = 1; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Thomas Hex") {
        $cell.Value2 = "Matthies Hex"
    }
}

# 5) Keep the worksheet's used range/dimension in sync.
$ws.Range("A1:T31").Select() | Out-Null

Write-Output "edit complete"
